$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D3").Value = 'id="Par1">The severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) has resulted in thousands of deaths in the world.

 Information about prediction model of prognosis of SARS-CoV-2 infection is scarce.

 We used machine learning for processing laboratory findings of 110 patients with SARS-CoV-2 pneumonia (including 51 non-survivors and 59 discharged patients).

 The maximum relevance minimum redundancy (mRMR) algorithm and the least absolute shrinkage and selection operator logistic regression model were used for selection of laboratory features.

 Seven laboratory features selected in the model were: prothrombin activity, urea, white blood cell, interleukin-2 receptor, indirect bilirubin, myoglobin, and fibrinogen degradation products.

 The signature constructed using the seven features had 98% [93%, 100%] sensitivity and 91% [84%, 99%] specificity in predicting outcome of SARS-CoV-2 pneumonia.

 Thus it is feasible to establish an accurate prediction model of outcome of SARS-CoV-2 pneumonia based on laboratory findings.

'
$ws.Range("E3").Value = '[Gang%Wu%NULL%1,  Shuchang%Zhou%NULL%1,  Yujin%Wang%NULL%1,  Wenzhi%Lv%NULL%2,  Shili%Wang%NULL%1,  Ting%Wang%751884926@qq.com%1,  Xiaoming%Li%lilyboston2002@qq.com%1]'
$ws.Range("D6").Value = 'Background
id="Par1">Coronavirus disease 2019 (COVID-19) has become a public health emergency of global concern.

 We aimed to explore the risk factors of 14-day and 28-day mortality and develop a model for predicting 14-day and 28-day survival probability among adult hospitalized patients with COVID-19.
Methods
id="Par2">In this multicenter, retrospective, cohort study, we examined 828 hospitalized patients with confirmed COVID-19 hospitalized in Wuhan Union Hospital and Central Hospital of Wuhan between January 12 and February 9, 2020. Among the 828 patients, 516 and 186 consecutive patients admitted in Wuhan Union Hospital were enrolled in the training cohort and the validation cohort, respectively.

 A total of 126 patients hospitalized in Central Hospital of Wuhan were enrolled in a second external validation cohort.

 Demographic, clinical, radiographic, and laboratory measures; treatment; proximate causes of death; and 14-day and 28-day mortality are described.

 Patients’ data were collected by reviewing the medical records, and their 14-day and 28-day outcomes were followed up.


Results
id="Par3">Of the 828 patients, 146 deaths were recorded until May 18, 2020. In the training set, multivariate Cox regression indicated that older age, lactate dehydrogenase level over 360 U/L, neutrophil-to-lymphocyte ratio higher than 8.0, and direct bilirubin higher than 5.0 μmol/L were independent predictors of 28-day mortality.

 Nomogram scoring systems for predicting the 14-day and 28-day survival probability of patients with COVID-19 were developed and exhibited strong discrimination and calibration power in the two external validation cohorts (C-index, 0.878 and 0.839).


Conclusion
id="Par4">Older age, high lactate dehydrogenase level, evaluated neutrophil-to-lymphocyte ratio, and high direct bilirubin level were independent predictors of 28-day mortality in adult hospitalized patients with confirmed COVID-19. The nomogram system based on the four factors revealed good discrimination and calibration, suggesting good clinical utility.


'
$ws.Range("E6").Value = '[Shuai%Zhang%NULL%1,  Mengfei%Guo%NULL%1,  Limin%Duan%NULL%1,  Feng%Wu%NULL%1,  Guorong%Hu%NULL%1,  Zhihui%Wang%NULL%1,  Qi%Huang%NULL%0,  Tingting%Liao%NULL%1,  Juanjuan%Xu%NULL%1,  Yanling%Ma%NULL%1,  Zhilei%Lv%NULL%1,  Wenjing%Xiao%NULL%1,  Zilin%Zhao%NULL%1,  Xueyun%Tan%NULL%1,  Daquan%Meng%NULL%1,  Shujing%Zhang%NULL%1,  E%Zhou%NULL%1,  Zhengrong%Yin%NULL%1,  Wei%Geng%NULL%1,  Xuan%Wang%NULL%1,  Jianchu%Zhang%NULL%1,  Jianguo%Chen%Chenj@mails.tjmu.edu.cn%1,  Yu%Zhang%whxhzy@163.com%0,  Yang%Jin%whuhjy@126.com%0]'
$ws.Range("D12").Value = 'Background
id="Par1">One-fifth of COVID-19 patients are seriously and critically ill cases and have a worse prognosis than non-severe cases.

 Although there is no specific treatment available for COVID-19, early recognition and supportive treatment may reduce the mortality.

 The aim of this study is to develop a functional nomogram that can be used by clinicians to estimate the risk of in-hospital mortality in patients hospitalized and treated for COVID-19 disease, and to compare the accuracy of model predictions with previous nomograms.


Methods
id="Par2">This retrospective study enrolled 709 patients who were over 18 years old and received inpatient treatment for COVID-19 disease.

 Multivariable Logistic Regression analysis was performed to assess the possible predictors of a fatal outcome.

 A nomogram was developed with the possible predictors and total point were calculated.


Results
id="Par3">Of the 709 patients treated for COVID-19, 75 (11%) died and 634 survived.

 The elder age, certain comorbidities (cancer, heart failure, chronic renal failure), dyspnea, lower levels of oxygen saturation and hematocrit, higher levels of C-reactive protein, aspartate aminotransferase and ferritin were independent risk factors for mortality.

 The prediction ability of total points was excellent (Area Under Curve = 0.922).


Conclusions
id="Par4">The nomogram developed in this study can be used by clinicians as a practical and effective tool in mortality risk estimation.

 So that with early diagnosis and intervention mortality in COVID-19 patients may be reduced.


'
$ws.Range("E12").Value = '[Hazal Cansu%Acar%hazal.acar@istanbul.edu.tr%1,  Günay%Can%NULL%2,  Günay%Can%NULL%0,  Rıdvan%Karaali%NULL%1,  Şermin%Börekçi%NULL%1,  İlker İnanç%Balkan%NULL%1,  Bilun%Gemicioğlu%NULL%1,  Dildar%Konukoğlu%NULL%1,  Ethem%Erginöz%NULL%1,  Mehmet Sarper%Erdoğan%NULL%1,  Fehmi%Tabak%NULL%1]'
$ws.Range("D14").Value = 'Objective
id="Par1">The aims of this study were to develop a multiparametric prognostic model for death in COVID-19 patients and to assess the incremental value of CT disease extension over clinical parameters.


Methods
id="Par2">Consecutive patients who presented to all five of the emergency rooms of the Reggio Emilia province between February 27 and March 23, 2020, for suspected COVID-19, underwent chest CT, and had a positive swab within 10 days were included in this retrospective study.

 Age, sex, comorbidities, days from symptom onset, and laboratory data were retrieved from institutional information systems.

 CT disease extension was visually graded as &lt; 20%, 20–39%, 40–59%, or ≥ 60%.

 The association between clinical and CT variables with death was estimated with univariable and multivariable Cox proportional hazards models; model performance was assessed using k-fold cross-validation for the area under the ROC curve (cvAUC).


Results
id="Par3">Of the 866 included patients (median age 59.8, women 39.2%), 93 (10.74%) died.

 Clinical variables significantly associated with death in multivariable model were age, male sex, HDL cholesterol, dementia, heart failure, vascular diseases, time from symptom onset, neutrophils, LDH, and oxygen saturation level.

 CT disease extension was also independently associated with death (HR = 7.56, 95% CI = 3.49; 16.38 for ≥ 60% extension).

 cvAUCs were 0.927 (bootstrap bias-corrected 95% CI = 0.899–0.947) for the clinical model and 0.936 (bootstrap bias-corrected 95% CI = 0.912–0.953) when adding CT extension.


Conclusions
id="Par4">A prognostic model based on clinical variables is highly accurate in predicting death in COVID-19 patients.

 Adding CT disease extension to the model scarcely improves its accuracy.


Key Points
id="Par5">
• Early identification of COVID-19 patients at higher risk of disease progression and death is crucial; the role of CT scan in defining prognosis is unclear.



Supplementary Information
The online version contains supplementary material available at 10.1007/s00330-021-07993-9.
'
$ws.Range("E14").Value = '[Giulia%Besutti%giulia.besutti@ausl.re.it%1,  Marta%Ottone%NULL%2,  Marta%Ottone%NULL%0,  Tommaso%Fasano%NULL%1,  Pierpaolo%Pattacini%NULL%1,  Valentina%Iotti%NULL%1,  Lucia%Spaggiari%NULL%1,  Riccardo%Bonacini%NULL%1,  Andrea%Nitrosi%NULL%1,  Efrem%Bonelli%NULL%1,  Simone%Canovi%NULL%1,  Rossana%Colla%NULL%1,  Alessandro%Zerbini%NULL%1,  Marco%Massari%NULL%1,  Ivana%Lattuada%NULL%1,  Anna Maria%Ferrari%NULL%1,  Paolo%Giorgi Rossi%NULL%1,  Massimo%Costantini%NULL%1,  Roberto%Grilli%NULL%1,  Massimiliano%Marino%NULL%1,  Giulio%Formoso%NULL%1,  Debora%Formisano%NULL%1,  Emanuela%Bedeschi%NULL%1,  Cinzia%Perilli%NULL%1,  Elisabetta%La Rosa%NULL%1,  Eufemia%Bisaccia%NULL%1,  Ivano%Venturi%NULL%1,  Massimo%Vicentini%NULL%1,  Cinzia%Campari%NULL%1,  Francesco%Gioia%NULL%1,  Serena%Broccoli%NULL%1,  Pamela%Mancuso%NULL%1,  Marco%Foracchia%NULL%1,  Mirco%Pinotti%NULL%1,  Nicola%Facciolongo%NULL%1,  Laura%Trabucco%NULL%1,  Stefano%De Pietri%NULL%1,  Giorgio Francesco%Danelli%NULL%1,  Laura%Albertazzi%NULL%1,  Enrica%Bellesia%NULL%1,  Mattia%Corradini%NULL%1,  Elena%Magnani%NULL%1,  Annalisa%Pilia%NULL%1,  Alessandra%Polese%NULL%1,  Silvia Storchi%Incerti%NULL%1,  Piera%Zaldini%NULL%1,  Bonanno%Orsola%NULL%1,  Matteo%Revelli%NULL%1,  Carlo%Salvarani%NULL%1,  Carmine%Pinto%NULL%1,  Francesco%Venturelli%NULL%1]'
$ws.Range("D16").Value = 'id="Par1">COVID-19 pandemic has created an extreme pressure on the global healthcare services.

 Fast, reliable, and early clinical assessment of the severity of the disease can help in allocating and prioritizing resources to reduce mortality.

 In order to study the important blood biomarkers for predicting disease mortality, a retrospective study was conducted on a dataset made public by Yan et al.

 in [1] of 375 COVID-19 positive patients admitted to Tongji Hospital (China) from January 10 to February 18, 2020. Demographic and clinical characteristics and patient outcomes were investigated using machine learning tools to identify key biomarkers to predict the mortality of individual patient.

 A nomogram was developed for predicting the mortality risk among COVID-19 patients.

 Lactate dehydrogenase, neutrophils (%), lymphocyte (%), high-sensitivity C-reactive protein, and age (LNLCA)—acquired at hospital admission—were identified as key predictors of death by multi-tree XGBoost model.

 The area under curve (AUC) of the nomogram for the derivation and validation cohort were 0.961 and 0.991, respectively.

 An integrated score (LNLCA) was calculated with the corresponding death probability.

 COVID-19 patients were divided into three subgroups: low-, moderate-, and high-risk groups using LNLCA cutoff values of 10.4 and 12.65 with the death probability less than 5%, 5–50%, and above 50%, respectively.

 The prognostic model, nomogram, and LNLCA score can help in early detection of high mortality risk of COVID-19 patients, which will help doctors to improve the management of patient stratification.

'
$ws.Range("E16").Value = '[Muhammad E. H.%Chowdhury%mchowdhury@qu.edu.qa%1,  Tawsifur%Rahman%NULL%2,  Tawsifur%Rahman%NULL%0,  Amith%Khandakar%NULL%3,  Somaya%Al-Madeed%NULL%2,  Susu M.%Zughaier%NULL%5,  Suhail A. R.%Doi%NULL%3,  Hanadi%Hassen%NULL%1,  Mohammad T.%Islam%NULL%1]'
$ws.Range("D33").Value = 'id="Par1">The epidemic phase of Coronavirus disease 2019 (COVID-19) made the Worldwide health system struggle against a severe interstitial pneumonia requiring high-intensity care settings for respiratory failure.

 A rationalisation of resources and a specific treatment path were necessary.

 The study suggests a predictive model drawing on clinical data gathered by 119 consecutive patients with laboratory-confirmed COVID-19 admitted in Busto Arsizio hospital.

 We derived a score that identifies the risk of clinical evolution and in-hospital mortality clustering patients into four groups.

 The study outcomes have been compared across the derivation and validation samples.

 The prediction rule is based on eight simple patient characteristics that were independently associated with study outcomes.

 It is able to stratify COVID-19 patients into four severity classes, with in-hospital mortality rates of 0% in group 1, 6–12.5% in group 2, 7–20% in group 3 and 60–86% in group 4 across the derivation and validation sample.

 The prediction model derived in this study identifies COVID-19 patients with low risk of in-hospital mortality and ICU admission.

 The prediction model that the study presents identifies COVID-19 patients with low risk of in-hospital mortality and admission to ICU.

 Moreover, it establishes an intermediate portion of patients that should be treated accurately in order to avoid an unfavourable clinical evolution.

 A further validation of the model is important before its implementation as a decision-making tool to guide the initial management of patients.

'
$ws.Range("E33").Value = '[Fabrizio%Foieni%fabrizio.foieni@asst-valleolona.it%1,  Girolamo%Sala%girolamo.sala@asst-valleolona.it%1,  Jason Giuseppe%Mognarelli%giuseppe.mognarelli@unimi.it%1,  Giulia%Suigo%NULL%1,  Davide%Zampini%NULL%1,  Matteo%Pistoia%NULL%1,  Mariella%Ciola%NULL%1,  Tommaso%Ciampani%NULL%1,  Carolina%Ultori%NULL%1,  Paolo%Ghiringhelli%NULL%1]'
$ws.Range("D48").Value = 'Background
id="Par1">Patients with severe Coronavirus Disease 2019 (COVID-19) will progress rapidly to acute respiratory failure or death.

 We aimed to develop a quantitative tool for early predicting mortality risk of patients with COVID-19.
Methods
id="Par2">301 patients with confirmed COVID-19 admitted to Main District and Tumor Center of the Union Hospital of Huazhong University of Science and Technology (Wuhan, China) between January 1, 2020 to February 15, 2020 were enrolled in this retrospective two-centers study.

 Data on patient demographic characteristics, laboratory findings and clinical outcomes was analyzed.

 A nomogram was constructed to predict the death probability of COVID-19 patients.


Results
id="Par3">Age, neutrophil-to-lymphocyte ratio, d-dimer and C-reactive protein obtained on admission were identified as predictors of mortality for COVID-19 patients by LASSO.

 The nomogram demonstrated good calibration and discrimination with the area under the curve (AUC) of 0.921 and 0.975 for the derivation and validation cohort, respectively.

 An integrated score (named ANDC) with its corresponding death probability was derived.

 Using ANDC cut-off values of 59 and 101, COVID-19 patients were classified into three subgroups.

 The death probability of low risk group (ANDC &lt; 59) was less than 5%, moderate risk group (59 ≤ ANDC ≤ 101) was 5% to 50%, and high risk group (ANDC &gt; 101) was more than 50%, respectively.


Conclusion
id="Par4">The prognostic nomogram exhibited good discrimination power in early identification of COVID-19 patients with high mortality risk, and ANDC score may help physicians to optimize patient stratification management.


'
$ws.Range("E48").Value = '[Zhihong%Weng%NULL%1,  Qiaosen%Chen%NULL%2,  Qiaosen%Chen%NULL%0,  Sumeng%Li%NULL%1,  Huadong%Li%NULL%1,  Qian%Zhang%NULL%1,  Sihong%Lu%NULL%1,  Li%Wu%NULL%1,  Leiqun%Xiong%NULL%1,  Bobin%Mi%NULL%1,  Di%Liu%NULL%1,  Mengji%Lu%NULL%1,  Dongliang%Yang%NULL%1,  Hongbo%Jiang%hongbojiang3@163.com%1,  Shaoping%Zheng%zhengspxx@126.com%1,  Xin%Zheng%xin11@hotmail.com%1]'
$ws.Range("D52").Value = 'Background
id="Par1">Accurately predicting patient outcomes in Severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) could aid patient management and allocation of healthcare resources.

 There are a variety of methods which can be used to develop prognostic models, ranging from logistic regression and survival analysis to more complex machine learning algorithms and deep learning.

 Despite several models having been created for SARS-CoV-2, most of these have been found to be highly susceptible to bias.

 We aimed to develop and compare two separate predictive models for death during admission with SARS-CoV-2.
Method
id="Par2">Between March 1 and April 24, 2020, 398 patients were identified with laboratory confirmed SARS-CoV-2 in a London teaching hospital.

 Data from electronic health records were extracted and used to create two predictive models using: (1) a Cox regression model and (2) an artificial neural network (ANN).

 Model performance profiles were assessed by validation, discrimination, and calibration.


Results
id="Par3">Both the Cox regression and ANN models achieved high accuracy (83.8%, 95% confidence interval (CI) 73.8–91.1 and 90.0%, 95% CI 81.2–95.6, respectively).

 The area under the receiver operator curve (AUROC) for the ANN (92.6%, 95% CI 91.1–94.1) was significantly greater than that of the Cox regression model (86.9%, 95% CI 85.7–88.2), p = 0.0136. Both models achieved acceptable calibration with Brier scores of 0.13 and 0.11 for the Cox model and ANN, respectively.


Conclusion
id="Par4">We demonstrate an ANN which is non-inferior to a Cox regression model but with potential for further development such that it can learn as new data becomes available.

 Deep learning techniques are particularly suited to complex datasets with non-linear solutions, which make them appropriate for use in conditions with a paucity of prior knowledge.

 Accurate prognostic models for SARS-CoV-2 can provide benefits at the patient, departmental and organisational level.


'
$ws.Range("E52").Value = '[Ahmed%Abdulaal%NULL%1,  Aatish%Patel%NULL%1,  Esmita%Charani%NULL%1,  Sarah%Denny%NULL%1,  Saleh A.%Alqahtani%NULL%1,  Gary W.%Davies%NULL%1,  Nabeela%Mughal%NULL%1,  Luke S. P.%Moore%l.moore@imperial.ac.uk%1]'
$ws.Range("D53").Value = 'id="Par1">The rapid spread of COVID-19 has resulted in the shortage of medical resources, which necessitates accurate prognosis prediction to triage patients effectively.

 This study used the nationwide cohort of South Korea to develop a machine learning model to predict prognosis based on sociodemographic and medical information.

 Of 10,237 COVID-19 patients, 228 (2.2%) died, 7772 (75.9%) recovered, and 2237 (21.9%) were still in isolation or being treated at the last follow-up (April 16, 2020).

 The Cox proportional hazards regression analysis revealed that age &gt; 70, male sex, moderate or severe disability, the presence of symptoms, nursing home residence, and comorbidities of diabetes mellitus (DM), chronic lung disease, or asthma were significantly associated with increased risk of mortality (p ≤ 0.047).

 For machine learning, the least absolute shrinkage and selection operator (LASSO), linear support vector machine (SVM), SVM with radial basis function kernel, random forest (RF), and k-nearest neighbors were tested.

 In prediction of mortality, LASSO and linear SVM demonstrated high sensitivities (90.7% [95% confidence interval: 83.3, 97.3] and 92.0% [85.9, 98.1], respectively) and specificities (91.4% [90.3, 92.5] and 91.8%, [90.7, 92.9], respectively) while maintaining high specificities &gt; 90%, as well as high area under the receiver operating characteristics curves (0.963 [0.946, 0.979] and 0.962 [0.945, 0.979], respectively).

 The most significant predictors for LASSO included old age and preexisting DM or cancer; for RF they were old age, infection route (cluster infection or infection from personal contact), and underlying hypertension.

 The proposed prediction model may be helpful for the quick triage of patients without having to wait for the results of additional tests such as laboratory or radiologic studies, during a pandemic when limited medical resources must be wisely allocated without hesitation.

'
$ws.Range("E53").Value = '[Chansik%An%NULL%1,  Hyunsun%Lim%NULL%1,  Dong-Wook%Kim%NULL%1,  Jung Hyun%Chang%NULL%1,  Yoon Jung%Choi%chris316@yuhs.ac%1,  Seong Woo%Kim%NULL%1]'
$ws.Range("D58").Value = 'id="Par1">Soaring cases of coronavirus disease (COVID-19) are pummeling the global health system.

 Overwhelmed health facilities have endeavored to mitigate the pandemic, but mortality of COVID-19 continues to increase.

 Here, we present a mortality risk prediction model for COVID-19 (MRPMC) that uses patients’ clinical data on admission to stratify patients by mortality risk, which enables prediction of physiological deterioration and death up to 20 days in advance.

 This ensemble model is built using four machine learning methods including Logistic Regression, Support Vector Machine, Gradient Boosted Decision Tree, and Neural Network.

 We validate MRPMC in an internal validation cohort and two external validation cohorts, where it achieves an AUC of 0.9621 (95% CI: 0.9464–0.9778), 0.9760 (0.9613–0.9906), and 0.9246 (0.8763–0.9729), respectively.

 This model enables expeditious and accurate mortality risk stratification of patients with COVID-19, and potentially facilitates more responsive health systems that are conducive to high risk COVID-19 patients.

'
$ws.Range("E58").Value = '[Yue%Gao%NULL%1,  Guang-Yao%Cai%NULL%1,  Wei%Fang%NULL%0,  Hua-Yi%Li%NULL%1,  Si-Yuan%Wang%NULL%2,  Si-Yuan%Wang%NULL%0,  Lingxi%Chen%NULL%1,  Yang%Yu%NULL%1,  Dan%Liu%NULL%0,  Sen%Xu%NULL%1,  Peng-Fei%Cui%NULL%1,  Shao-Qing%Zeng%NULL%2,  Shao-Qing%Zeng%NULL%0,  Xin-Xia%Feng%NULL%1,  Rui-Di%Yu%NULL%1,  Ya%Wang%NULL%1,  Yuan%Yuan%NULL%1,  Xiao-Fei%Jiao%NULL%1,  Jian-Hua%Chi%NULL%1,  Jia-Hao%Liu%NULL%1,  Ru-Yuan%Li%NULL%1,  Xu%Zheng%NULL%1,  Chun-Yan%Song%NULL%1,  Ning%Jin%NULL%1,  Wen-Jian%Gong%NULL%1,  Xing-Yu%Liu%NULL%1,  Lei%Huang%NULL%0,  Xun%Tian%NULL%1,  Lin%Li%NULL%1,  Hui%Xing%NULL%1,  Ding%Ma%NULL%1,  Chun-Rui%Li%NULL%1,  Fei%Ye%yeyuanbei@hotmail.com%1,  Qing-Lei%Gao%qingleigao@hotmail.com%2,  Qing-Lei%Gao%qingleigao@hotmail.com%0]'
$ws.Range("D61").Value = 'Background
id="Par1">COVID-19 has become a major global threat.

 The present study aimed to develop a nomogram model to predict the survival of COVID-19 patients based on their clinical and laboratory data at admission.


Methods
id="Par2">COVID-19 patients who were admitted at Hankou Hospital and Huoshenshan Hospital in Wuhan, China from January 12, 2020 to March 20, 2020, whose outcome during the hospitalization was known, were retrospectively reviewed.

 The categorical variables were compared using Pearson’s χ2-test or Fisher’s exact test, and continuous variables were analyzed using Student’s t-test or Mann Whitney U-test, as appropriate.

 Then, variables with a P-value of ≤0.1 were included in the log-binomial model, and merely these independent risk factors were used to establish the nomogram model.

 The discrimination of the nomogram was evaluated using the area under the receiver operating characteristic curve (AUC), and internally verified using the Bootstrap method.


Results
id="Par3">A total of 262 patients (134 surviving and 128 non-surviving patients) were included in the analysis.

 Seven variables, which included age (relative risk [RR]: 0.905, 95% confidence interval [CI]: 0.868–0.944; P &lt; 0.001), chronic heart disease (CHD, RR: 0.045, 95% CI: 0.0097–0.205; P &lt; 0.001, the percentage of lymphocytes (Lym%, RR: 1.125, 95% CI: 1.041–1.216; P = 0.0029), platelets (RR: 1.008, 95% CI: 1.003–1.012; P = 0.001), C-reaction protein (RR: 0.982, 95% CI: 0.973–0.991; P &lt; 0.001), lactate dehydrogenase (LDH, RR: 0.993, 95% CI: 0.990–0.997; P &lt; 0.001) and D-dimer (RR: 0.734, 95% CI: 0.617–0.879; P &lt; 0.001), were identified as the independent risk factors.

 The nomogram model based on these factors exhibited a good discrimination, with an AUC of 0.948 (95% CI: 0.923–0.973).


Conclusions
id="Par4">A nomogram based on age, CHD, Lym%, platelets, C-reaction protein, LDH and D-dimer was established to accurately predict the prognosis of COVID-19 patients.

 This can be used as an alerting tool for clinicians to take early intervention measures, when necessary.


'
$ws.Range("E61").Value = '[Xiaojun%Ma%NULL%1,  Huifang%Wang%NULL%1,  Junwei%Huang%NULL%1,  Yan%Geng%NULL%1,  Shuqi%Jiang%NULL%1,  Qiuping%Zhou%NULL%1,  Xuan%Chen%NULL%1,  Hongping%Hu%NULL%1,  Weifeng%Li%NULL%1,  Chengbin%Zhou%NULL%1,  Xinglin%Gao%NULL%1,  Na%Peng%pnatz@163.com%1,  Yiyu%Deng%yiyudeng666@163.com%1]'
